$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet "DA_Price" -> "Price" (for naming consistency)
$ws.Name = "Price"

# Move the active selection from F14 to E12
$ws.Range("E12").Select() | Out-Null
